# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet between "2021-Q4" and "总计", fills it
# with the fund-holding detail rows for the new quarter, and refreshes the
# "总计" (totals) sheet with a new leading row for 2022-Q1 (pushing the
# existing 2021-Q4 total row down).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$wsTotalOld = $wb.Worksheets.Item("总计")

# Recreate "总计" after inserting the new quarter sheet so tab order /
# sheetId allocation ends up as [2021-Q4, 2022-Q1, 总计].
$wsTotalOld.Delete()

# ---------------------------------------------------------------------
# New sheet: 2022-Q1
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add($null, $ws1)
$wsQ1.Name = "2022-Q1"

$headerRangeQ1 = $wsQ1.Range("B1:H1")
$headerRangeQ1.Font.Bold = $true
$headerRangeQ1.HorizontalAlignment = -4108
$headerRangeQ1.VerticalAlignment = -4160
$headerRangeQ1.Borders.LineStyle = 1

$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

$idxCellQ1 = $wsQ1.Range("A2")
$idxCellQ1.Font.Bold = $true
$idxCellQ1.HorizontalAlignment = -4108
$idxCellQ1.VerticalAlignment = -4160
$idxCellQ1.Borders.LineStyle = 1
$idxCellQ1.Value = 0

$wsQ1.Range("B2").Value = "'202801"
$wsQ1.Range("C2").Value = "南方全球精选配置(QDII-FOF)"
$wsQ1.Range("D2").Value = "'18.00"
$wsQ1.Range("E2").Value = "'28.82"
$wsQ1.Range("F2").Value = "'1.11"
$wsQ1.Range("G2").Value = "'0.1998"
$wsQ1.Range("H2").Value = 7

# ---------------------------------------------------------------------
# New sheet: 总计 (rebuilt after Q1, with a fresh leading row for 2022-Q1)
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

$headerRangeTotal = $wsTotal.Range("B1:D1")
$headerRangeTotal.Font.Bold = $true
$headerRangeTotal.HorizontalAlignment = -4108
$headerRangeTotal.VerticalAlignment = -4160
$headerRangeTotal.Borders.LineStyle = 1

$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$idxCell2 = $wsTotal.Range("A2")
$idxCell2.Font.Bold = $true
$idxCell2.HorizontalAlignment = -4108
$idxCell2.VerticalAlignment = -4160
$idxCell2.Borders.LineStyle = 1
$idxCell2.Value = 0

$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.2

$idxCell3 = $wsTotal.Range("A3")
$idxCell3.Font.Bold = $true
$idxCell3.HorizontalAlignment = -4108
$idxCell3.VerticalAlignment = -4160
$idxCell3.Borders.LineStyle = 1
$idxCell3.Value = 1

$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 6
$wsTotal.Range("D3").Value = 0.58

# Restore original active sheet (2021-Q4).
$ws1.Activate()
